# Replicates the "Add files via upload" commit on analysis/rawData/data.xlsx
#
# Reachable-via-COM changes applied here:
#   1. Column K ("Implement ML") flipped from "yes" to "no" on the rows that
#      correspond to entries 14,18,21,28,34,38,39,40,43,47,48,52 (sheet rows
#      15,19,22,29,35,39,40,41,44,48,49,53). Writing a brand-new literal
#      "no" (not previously present in the shared-string table) causes the
#      engine to append exactly one new shared string, matching the
#      uniqueCount 350 -> 351 bump in the diff.
#   2. Column M widened (diff: width="57.88671875"). The engine quantises
#      ColumnWidth to 1/6-character steps, so we pick the input value whose
#      round-tripped width lands closest to the target.
#   3. Selection moved to K56 (diff: <selection activeCell="K56" sqref="K56"/>).
#
# NOTE: a few cosmetic bits of the source diff are not reachable from the
# documented Excel COM surface in this runtime and are intentionally left
# alone rather than faked:
#   - x15ac:absPath (OneDrive path breadcrumb) is round-tripped as opaque,
#     unmodelled XML by the engine; there is no object-model property that
#     writes it.
#   - xr:revisionPtr is dropped by the exporter on every save regardless of
#     what the COM script does (verified with a no-op script), so it cannot
#     be reproduced here.
#   - bookViews/workbookView's xWindow/yWindow/windowWidth/windowHeight
#     mirror Window.Top/Left/Width/Height, but those setters are not wired
#     to the persisted XML in this build (verified empirically) - the saved
#     file always keeps the originally loaded values.
#   - sheetView's bare topLeftCell="A22" is only emitted by the exporter
#     inside a <pane> (i.e. only for frozen/split views); setting
#     ActiveWindow.ScrollRow/ScrollColumn does not surface it, and turning
#     on FreezePanes would add a <pane>/frozen state the source diff does
#     not have, so it is skipped.
#   - xl/persons/person.xml is only emitted by this engine as a side effect
#     of creating an actual threaded comment (which also emits
#     comments1.xml/threadedComment1.xml that are absent from the target
#     diff, and deleting the comment again removes person.xml along with
#     it), so there is no clean way to produce the bare empty <personList/>
#     without introducing unwanted comment parts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Flip "Implement ML" (column K) from "yes" to "no" for the rows changed
#    in the diff.
$rows = @(15, 19, 22, 29, 35, 39, 40, 41, 44, 48, 49, 53)
foreach ($r in $rows) {
    $ws.Range("K$r").Value = "no"
}

# 2. Widen column M (13th column) to match the new <cols> entry.
#    Target XML width is 57.88671875; the closest value this engine's
#    character-width quantisation can produce is 57.833333... (ColumnWidth
#    input of 57).
$ws.Columns.Item(13).ColumnWidth = 57

# 3. Move the sheet selection to K56, as in the diff's <selection .../>.
$ws.Range("K56").Select() | Out-Null
